# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, pushing the existing "Late" / "heading" / "Outstanding" columns one
# place to the right, then re-select that sheet with a new active cell.

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N (shifts N->O, O->P, P->Q)
$wsRepayment.Columns("N").Insert()

# Give the newly inserted column the same width as column M (10.7109375
# characters) instead of leaving it at the default sheet width.
$wsRepayment.Columns("N").ColumnWidth = $wsRepayment.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet and move the selection to S8,
# which becomes the new active tab / selection stored in the workbook.
$wsRepayment.Activate()
$wsRepayment.Range("S8").Select()
